# Primeiras verificações de entradas
# Adds the first batch of delivery-check rows (rows 3-8) to the "Dados" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain TEXT (so date/time-looking strings such as
# "02/05/2005" or "00:00" are not auto-converted into date/time serials).
# A leading apostrophe is the classic "force text" entry method, after which
# the "Normal" style is restored so no residual number-format/style index is
# left behind on the cell.
function Set-TextValue($row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

# Helper: write a genuine numeric value.
function Set-NumberValue($row, $col, $val) {
    $ws.Cells.Item($row, $col).Value = $val
}

$rows = @(
    @{ Row=3; A=2; B="JEan";   C="02/05/2005";  D="00:00"; E=0; F=0; G=0;  H=0;  K="Pendente" },
    @{ Row=4; A=3; B="Jean";   C="32/13/20221"; D="10:20"; E=0; F=0; G=0;  H=0;  K="Pendente" },
    @{ Row=5; A=4; B="Jean";   C="01/01/2021";  D="10:33"; E=0; F=0; G=0;  H=0;  K="Pendente" },
    @{ Row=6; A=5; B="AFSDFS"; C="12/13/2022";  D="10:45"; E=0; F=0; G=0;  H=0;  K="Pendente" },
    @{ Row=7; A=6; B="Jean";   C="01/01/2022";  D="23:52"; E=1; F=0; G=0;  H=0;  K="Pendente" },
    @{ Row=8; A=7; B="Jean";   C="31/01/2022";  D="10:50"; E=1; F=0; G=10; H=15; K="Pendente" }
)

foreach ($r in $rows) {
    Set-NumberValue $r.Row 1 $r.A
    Set-TextValue   $r.Row 2 $r.B
    Set-TextValue   $r.Row 3 $r.C
    Set-TextValue   $r.Row 4 $r.D
    Set-NumberValue $r.Row 5 $r.E
    Set-NumberValue $r.Row 6 $r.F
    Set-NumberValue $r.Row 7 $r.G
    Set-NumberValue $r.Row 8 $r.H
    Set-TextValue   $r.Row 11 $r.K
}

# Row 8 also carries an (empty) "Mensagem adicional" entry in column J — a
# present-but-blank text cell. A lone apostrophe forces Excel to store an
# empty text value instead of leaving the cell completely absent.
Set-TextValue 8 10 ""
